$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: remove old row 2 ("Fonte Carregador ...") entirely; rows below shift up
$ws.Rows.Item(2).Delete()

# Step 2: insert two new columns before column A, shifting nome..link from A-G to C-I
$ws.Range("A1:B1").EntireColumn.Insert()

# Step 3: give the new header cells (A1, B1) the same bold/border/centered style as the other headers
$ws.Cells.Item(1,3).Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Step 4: header captions for the new columns
$ws.Cells.Item(1,1).Value = "data"
$ws.Cells.Item(1,2).Value = "loja"

# Step 5: write full data block (rows 2-9, columns A-I) matching the target content/order
# row 2: 'Controle Longa Distancia Jfa Acqua 1200 Resistente'
$ws.Cells.Item(2,1).Value = "30/07/2024"
$ws.Cells.Item(2,2).Value = "vetaacessorios"
$ws.Cells.Item(2,3).Value = "Controle Longa Distancia Jfa Acqua 1200 Resistente A Água"
$ws.Cells.Item(2,4).Value = "ACQUA"
$ws.Cells.Item(2,5).Value = 75
$ws.Cells.Item(2,6).Value = "Baixo"
$ws.Cells.Item(2,7).Value = "NA"
$ws.Cells.Item(2,8).Value = "classico"
$ws.Cells.Item(2,9).Value = "https://www.mercadolivre.com.br/controle-longa-distancia-jfa-acqua-1200-resistente-a-agua/p/MLB28961390?pdp_filters=seller_id:1162748365#searchVariation=MLB28961390&position=1&search_layout=stack&type=product&tracking_id=f0e739c8-da54-412c-bb96-855c84702b37"

# row 3: 'Controle Jfa Acqua K1200 Longa Distancia Completo '
$ws.Cells.Item(3,1).Value = "30/07/2024"
$ws.Cells.Item(3,2).Value = "vetaacessorios"
$ws.Cells.Item(3,3).Value = "Controle Jfa Acqua K1200 Longa Distancia Completo Top"
$ws.Cells.Item(3,4).Value = "ACQUA"
$ws.Cells.Item(3,5).Value = 67.5
$ws.Cells.Item(3,6).Value = "Baixo"
$ws.Cells.Item(3,7).Value = "NA"
$ws.Cells.Item(3,8).Value = "classico"
$ws.Cells.Item(3,9).Value = "https://produto.mercadolivre.com.br/MLB-3037029276-controle-jfa-acqua-k1200-longa-distancia-completo-top-_JM#position%3D3%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df0e739c8-da54-412c-bb96-855c84702b37"

# row 4: 'Controle Longa Distância Jfa Redline 1200 Metros V'
$ws.Cells.Item(4,1).Value = "30/07/2024"
$ws.Cells.Item(4,2).Value = "vetaacessorios"
$ws.Cells.Item(4,3).Value = "Controle Longa Distância Jfa Redline 1200 Metros Vermelho"
$ws.Cells.Item(4,4).Value = "Sem Modelo"
$ws.Cells.Item(4,5).Value = 78.9
$ws.Cells.Item(4,6).Value = ""
$ws.Cells.Item(4,7).Value = "NA"
$ws.Cells.Item(4,8).Value = "classico"
$ws.Cells.Item(4,9).Value = "https://produto.mercadolivre.com.br/MLB-3860722412-controle-longa-distncia-jfa-redline-1200-metros-vermelho-_JM#position%3D4%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df0e739c8-da54-412c-bb96-855c84702b37"

# row 5: 'Controle Longa Distancia Jfa K1200 Acqua Completo '
$ws.Cells.Item(5,1).Value = "30/07/2024"
$ws.Cells.Item(5,2).Value = "vetaacessorios"
$ws.Cells.Item(5,3).Value = "Controle Longa Distancia Jfa K1200 Acqua Completo Top Preto"
$ws.Cells.Item(5,4).Value = "ACQUA"
$ws.Cells.Item(5,5).Value = 77.9
$ws.Cells.Item(5,6).Value = "Baixo"
$ws.Cells.Item(5,7).Value = "NA"
$ws.Cells.Item(5,8).Value = "classico"
$ws.Cells.Item(5,9).Value = "https://produto.mercadolivre.com.br/MLB-2927266757-controle-longa-distancia-jfa-k1200-acqua-completo-top-preto-_JM#position%3D5%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df0e739c8-da54-412c-bb96-855c84702b37"

# row 6: 'Controle Remoto Jfa Longa Distancia K1200 Acqua Co'
$ws.Cells.Item(6,1).Value = "30/07/2024"
$ws.Cells.Item(6,2).Value = "vetaacessorios"
$ws.Cells.Item(6,3).Value = "Controle Remoto Jfa Longa Distancia K1200 Acqua Completo Top"
$ws.Cells.Item(6,4).Value = "ACQUA"
$ws.Cells.Item(6,5).Value = 75
$ws.Cells.Item(6,6).Value = "Baixo"
$ws.Cells.Item(6,7).Value = "NA"
$ws.Cells.Item(6,8).Value = "classico"
$ws.Cells.Item(6,9).Value = "https://produto.mercadolivre.com.br/MLB-3037065409-controle-remoto-jfa-longa-distancia-k1200-acqua-completo-top-_JM#position%3D6%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df0e739c8-da54-412c-bb96-855c84702b37"

# row 7: 'Controle Longa Distancia Jfa Acqua  1200 Mts  Pret'
$ws.Cells.Item(7,1).Value = "30/07/2024"
$ws.Cells.Item(7,2).Value = "vetaacessorios"
$ws.Cells.Item(7,3).Value = "Controle Longa Distancia Jfa Acqua  1200 Mts  Preto Completo"
$ws.Cells.Item(7,4).Value = "ACQUA"
$ws.Cells.Item(7,5).Value = 109.9
$ws.Cells.Item(7,6).Value = "Acima"
$ws.Cells.Item(7,7).Value = "NA"
$ws.Cells.Item(7,8).Value = "classico"
$ws.Cells.Item(7,9).Value = "https://produto.mercadolivre.com.br/MLB-3037013938-controle-longa-distancia-jfa-acqua-1200-mts-preto-completo-_JM#position%3D7%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df0e739c8-da54-412c-bb96-855c84702b37"

# row 8: 'Controle Longa Distância Jfa Redline 1200 Metros V'
$ws.Cells.Item(8,1).Value = "30/07/2024"
$ws.Cells.Item(8,2).Value = "vetaacessorios"
$ws.Cells.Item(8,3).Value = "Controle Longa Distância Jfa Redline 1200 Metros Vermelho"
$ws.Cells.Item(8,4).Value = "Sem Modelo"
$ws.Cells.Item(8,5).Value = 110
$ws.Cells.Item(8,6).Value = ""
$ws.Cells.Item(8,7).Value = "NA"
$ws.Cells.Item(8,8).Value = "classico"
$ws.Cells.Item(8,9).Value = "https://produto.mercadolivre.com.br/MLB-2731131087-controle-longa-distncia-jfa-redline-1200-metros-vermelho-_JM#position%3D8%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df0e739c8-da54-412c-bb96-855c84702b37"

# row 9: 'Controle Longa Distância Jfa Redline 1200 Metros V'
$ws.Cells.Item(9,1).Value = "30/07/2024"
$ws.Cells.Item(9,2).Value = "vetaacessorios"
$ws.Cells.Item(9,3).Value = "Controle Longa Distância Jfa Redline 1200 Metros Vermelho"
$ws.Cells.Item(9,4).Value = "Sem Modelo"
$ws.Cells.Item(9,5).Value = 78.9
$ws.Cells.Item(9,6).Value = ""
$ws.Cells.Item(9,7).Value = "NA"
$ws.Cells.Item(9,8).Value = "classico"
$ws.Cells.Item(9,9).Value = "https://produto.mercadolivre.com.br/MLB-3860735858-controle-longa-distncia-jfa-redline-1200-metros-vermelho-_JM#position%3D9%26search_layout%3Dstack%26type%3Ditem%26tracking_id%3Df0e739c8-da54-412c-bb96-855c84702b37"
